$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 16, pushing existing rows 16-34 down to 19-37.
$ws.Rows("16:18").Insert()

# Row 16: Chirimoya - Especial, week of 2021-09-21, Provincia de Limari
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44460
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107002
$ws.Range("J16").Value = "Chirimoya"
$ws.Range("K16").Value = "Cultivar IV Región"
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 2600
$ws.Range("O16").Value = 2700
$ws.Range("P16").Value = 2650
$ws.Range("Q16").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 2650
$ws.Range("T16").Value = 1

# Row 17: Chirimoya - Primera, week of 2021-09-21, Provincia de Limari
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Terminal La Palmera de La Serena"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = 44460
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = "Chirimoya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 240
$ws.Range("N17").Value = 2200
$ws.Range("O17").Value = 2300
$ws.Range("P17").Value = 2250
$ws.Range("Q17").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2250
$ws.Range("T17").Value = 1

# Row 18: Chirimoya - Segunda, week of 2021-09-21, Provincia de Limari
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "Terminal La Palmera de La Serena"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44460
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = "Otros"
$ws.Range("I18").Value = 100107002
$ws.Range("J18").Value = "Chirimoya"
$ws.Range("K18").Value = "Cultivar IV Región"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 1900
$ws.Range("O18").Value = 2000
$ws.Range("P18").Value = 1950
$ws.Range("Q18").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R18").Value = "Provincia de Limarí"
$ws.Range("S18").Value = 1950
$ws.Range("T18").Value = 1
